$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.685.09'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = '  +0.61%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.007.32'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = '  +3.13%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '381.34'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '106.45'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +1.68%  '
$ws.Range("E7").Value = '  +0.89%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.602'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  +1.90%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.93'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +2.71%  '
$ws.Range("E11").Value = '  +0.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0849'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  +1.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.89'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  +1.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.483.43'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = '  +2.97%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.55'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  +2.66%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.005.74'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  +2.58%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.977'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  +2.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.720.46'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  +0.91%  '
$ws.Range("E19").Value = '  +4.88%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.47'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  +3.19%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.16'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +0.73%  '
$ws.Range("E22").Value = '  +2.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.97'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  +0.76%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '264.47'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  +1.51%  '
$ws.Range("E25").Value = '  +3.76%  '
$ws.Range("E26").Value = '  -0.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.33'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  +18.88%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.53'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  +3.51%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '26.24'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  +0.39%  '
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("E31").Value = '  -1.83%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.98'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  -0.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.99'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  -0.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '51.22'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +0.23%  '
$ws.Range("E35").Value = '  -3.67%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0448'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  +5.59%  '
$ws.Range("E37").Value = '  -0.09%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.12'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  -1.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '17.65'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +3.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.64'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  -6.25%  '
$ws.Range("E41").Value = '  +0.49%  '
$ws.Range("E42").Value = '  +2.94%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '122.69'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  +2.31%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.46'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  -0.05%  '
$ws.Range("E45").Value = '  -2.87%  '
$ws.Range("E46").Value = '  +7.38%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.275'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  +15.57%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.061.22'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  -0.89%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.30'
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = '  +2.95%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0355'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +13.22%  '
$ws.Range("E51").Value = '  +3.20%  '
